# Generate Report for Handback
# Refresh the handoff/handback timestamps for the "30ea3542-..." row after a
# new handback report run. Only the zh-cn and de-de rows for that source
# file get new "Correspond Handoff/Handback DateTime" stamps, and the
# Overview sheet's "Latest HO Xliff Generate Date" for that same file.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_overview.Range("G2").Value = "2016-09-03 22:51:52"

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("H2").Value = "2016-09-03 22:51:48"
$ws_zhcn.Range("K2").Value = "2016-09-03 22:52:09"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("H2").Value = "2016-09-03 22:51:52"
$ws_dede.Range("K2").Value = "2016-09-03 22:52:16"
